$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Importe") holds amounts stored as text (Argentine-locale formatted,
# e.g. "144.000,00"). A scraping bug mangled the locale conversion, so the fix
# re-writes each affected amount as plain "144000.00" text (no thousands separator,
# period as decimal separator). Force Text number format first so Excel keeps the
# value as a string instead of re-parsing it into a real number.
$amountCells = @{
    "H2" = "144000.00"
    "H3" = "54000.00"
    "H4" = "53000.00"
    "H5" = "54000.00"
    "H6" = "350500.00"
    "H7" = "1756000.00"
    "H8" = "20.00"
    "H9" = "193303.13"
    "H10" = "162140.00"
    "H11" = "150000.00"
    "H12" = "32.97"
    "H13" = "2202.20"
    "H14" = "103111.75"
    "H15" = "440354.57"
    "H16" = "12133.00"
    "H17" = "7378.00"
    "H18" = "168425.08"
    "H19" = "17315.00"
    "H20" = "170.00"
    "H21" = "662.50"
    "H22" = "2058.00"
    "H23" = "613.30"
    "H24" = "26671.56"
    "H25" = "2376.00"
    "H26" = "9774.06"
    "H27" = "6376.51"
    "H28" = "2000.00"
    "H29" = "233.30"
    "H30" = "30.00"
    "H31" = "6370.00"
    "H32" = "6210.00"
    "H33" = "1163.73"
    "H34" = "398.84"
    "H35" = "95.00"
    "H36" = "212.00"
    "H37" = "134932.50"
    "H38" = "81328.09"
    "H39" = "2373.92"
    "H40" = "22.12"
    "H41" = "2167.90"
    "H42" = "236.50"
    "H43" = "14274.80"
    "H44" = "55512.81"
    "H45" = "143.50"
    "H46" = "140.84"
    "H47" = "28876.24"
    "H48" = "2566.32"
    "H49" = "575.00"
    "H50" = "7958.00"
    "H51" = "14991.82"
    "H52" = "1620.00"
    "H53" = "638.00"
    "H54" = "218.40"
    "H55" = "49036.27"
    "H56" = "75.00"
    "H57" = "6573.14"
    "H58" = "42.00"
    "H59" = "8811.10"
    "H60" = "2018.00"
    "H61" = "138500.00"
    "H62" = "19900.00"
    "H63" = "13158.00"
    "H64" = "650.00"
    "H65" = "694.20"
    "H66" = "3.81"
    "H67" = "122754.30"
    "H68" = "14.26"
    "H69" = "19360.34"
    "H70" = "175613.40"
    "H71" = "1015.26"
    "H72" = "417.48"
    "H73" = "6738.50"
    "H74" = "782.38"
    "H75" = "1515.00"
    "H76" = "9900.00"
    "H77" = "50.62"
    "H78" = "334.00"
    "H79" = "83.76"
    "H80" = "2114.81"
    "H81" = "8800.00"
    "H82" = "250.00"
    "H83" = "670.00"
    "H84" = "14.80"
    "H85" = "35200.00"
    "H86" = "16176.00"
    "H87" = "8590.00"
    "H88" = "16495.00"
    "H89" = "6510.00"
    "H90" = "1000.00"
    "H91" = "532.00"
    "H92" = "340.00"
    "H93" = "12647.00"
    "H94" = "3140.00"
    "H95" = "2071.00"
    "H96" = "536.00"
    "H97" = "88.00"
    "H98" = "3540.00"
    "H99" = "63.76"
    "H100" = "10411.03"
    "H101" = "120.00"
    "H102" = "209.96"
    "H103" = "289935.82"
    "H104" = "18945.38"
    "H105" = "329.32"
    "H106" = "566.00"
    "H107" = "36000.00"
    "H108" = "11241.90"
    "H109" = "26943.85"
    "H110" = "6462.00"
    "H111" = "604.00"
    "H112" = "50386.80"
    "H113" = "1602.25"
    "H114" = "3143.00"
    "H115" = "2160.00"
    "H116" = "29052.32"
    "H117" = "2957.00"
    "H118" = "2808.47"
    "H119" = "23.00"
    "H120" = "1194.00"
    "H121" = "200.00"
    "H122" = "4889.17"
    "H123" = "5268.00"
    "H124" = "880.00"
    "H125" = "1400.00"
    "H126" = "291.20"
    "H127" = "58.94"
    "H128" = "5412.00"
    "H129" = "9408.00"
    "H130" = "3200.00"
    "H131" = "910.00"
    "H132" = "17300.00"
    "H133" = "376.00"
    "H134" = "1312.00"
    "H135" = "1500.00"
    "H136" = "6000.00"
    "H137" = "3200.00"
    "H138" = "59.90"
    "H139" = "96.00"
    "H140" = "126.50"
    "H141" = "93012.00"
    "H142" = "580.54"
    "H143" = "685999.82"
    "H144" = "7985.76"
    "H145" = "1000.00"
    "H146" = "5616.00"
    "H147" = "6700.00"
    "H148" = "6050.00"
    "H149" = "5200.00"
    "H150" = "1300.00"
    "H151" = "750.00"
    "H152" = "565.00"
    "H153" = "250.00"
    "H154" = "20.00"
    "H155" = "833.66"
    "H156" = "198.00"
    "H157" = "5600.00"
    "H158" = "364.60"
    "H159" = "5925.00"
    "H160" = "4820.00"
    "H161" = "294.36"
    "H162" = "30.00"
    "H163" = "6900.00"
    "H164" = "275.54"
    "H165" = "900.00"
    "H166" = "1180.00"
    "H167" = "6713.56"
    "H168" = "4277.66"
    "H169" = "3761.00"
    "H170" = "1716.00"
    "H171" = "2157.60"
    "H172" = "3479.02"
    "H173" = "3720.00"
    "H174" = "573.00"
    "H175" = "3100.00"
    "H176" = "15705.87"
    "H177" = "2629.00"
    "H178" = "2950.00"
    "H179" = "3352.20"
    "H180" = "60.00"
    "H181" = "600.00"
    "H182" = "20500.00"
    "H183" = "26038.80"
    "H184" = "229.65"
    "H185" = "43500.00"
    "H186" = "4582.51"
    "H187" = "26400.00"
    "H188" = "2849.52"
    "H189" = "819951.34"
    "H190" = "14500.00"
    "H191" = "4613.60"
    "H192" = "179836.00"
    "H193" = "336780.00"
    "H194" = "55000.00"
    "H195" = "217500.00"
    "H196" = "75000.00"
    "H197" = "135166.00"
    "H198" = "244344.00"
    "H199" = "68420.00"
    "H200" = "137948.00"
    "H201" = "268304.00"
    "H202" = "154800.00"
    "H203" = "13110.00"
    "H204" = "4500.00"
    "H205" = "29800.00"
    "H206" = "240000.00"
}
foreach ($addr in $amountCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $amountCells[$addr]
}

# A handful of "Razon social" / "Nombre Fantasia" entries used ", " to separate
# co-owners (e.g. "FERNANDEZ, MARIO HUGO"). The same buggy locale-fix pass
# stripped periods and turned commas into periods there too, so replicate that.
$nameCells = @{
    "E44" = "RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH"
    "E56" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E92" = "FERNANDEZ. MARIO HUGO"
    "E94" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E97" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "F97" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "E115" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "F115" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "E124" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
    "F136" = "MERCANZINI. GASTON ARIEL"
    "E156" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E173" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
}
foreach ($addr in $nameCells.Keys) {
    $ws.Range($addr).Value = $nameCells[$addr]
}

Write-Output "done"
